# The single Bibliografia run is split into multiple runs, one per
# reference entry, joined by manual line breaks (<w:br/>), with a
# trailing space (xml:space="preserve") kept on entries that had one.
$d = $word.ActiveDocument

$startMarker = "INCROPERA, F, P"
$endMarker   = "Academic Press, 1979."

$fullText = $d.Content.Text
$start = $fullText.IndexOf($startMarker)
$end   = $fullText.IndexOf($endMarker) + $endMarker.Length

$target = $d.Range($start, $end)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>' + `
  '<w:t xml:space="preserve">INCROPERA, F, P; DEWITT, D. P. Fundamentos de Transferência de Calor e de Massa, LTC Editora, 2005. </w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">BENNETT, C. D.; MYERS, J. E. Fenômenos de Transporte. McGraw-Hill. </w:t>' + `
  '<w:br/>' + `
  '<w:t>KREITH, F.; BOHN, M. S. Princípios de Transferência de Calor, Thomson Learning, 2003.</w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004. </w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">FOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001. </w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">SISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988. </w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">HOLMAN, J. P. Transferência de Calor, McGraw-Hill, 1983. </w:t>' + `
  '<w:br/>' + `
  '<w:t>POIRIER, D.R.; GEIGER, G.H. Transport Phenomena in Materials Processing, TMS, 1994.</w:t>' + `
  '<w:br/>' + `
  '<w:t xml:space="preserve">GASKELL, David R. Introduction to Transport Phenomena in Materials Engineering. Prentice Hall, 1991. </w:t>' + `
  '<w:br/>' + `
  '<w:t>SZEKELY, J. Fluid Flow Phenomena in Metals Processing. Academic Press, 1979.</w:t>' + `
  '</w:r></w:p>'

$target.InsertXML($newXml)
Write-Output "Done"
